# Adds a new "2022-Q1" sheet (fund-holdings detail) positioned right before
# the existing "总计" (summary) sheet, and inserts a corresponding new
# top data-row ("2022-Q1", 26, 2.64) into the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Insert the new "2022-Q1" worksheet between "2021-Q4" and "总计".
# ---------------------------------------------------------------------
$anchor = $wb.Worksheets.Item("2021-Q4")
$ws = $wb.Worksheets.Add($null, $anchor)
$ws.Name = "2022-Q1"

# Header row (bold/centered style already lives on the sheet via column A's
# style elsewhere; here we just set the header text as the source sheets do).
$ws.Cells.Item(1, 2).Value = "基金代码"
$ws.Cells.Item(1, 3).Value = "基金名称"
$ws.Cells.Item(1, 4).Value = "基金规模"
$ws.Cells.Item(1, 5).Value = "股票总仓位"
$ws.Cells.Item(1, 6).Value = "仓位占比"
$ws.Cells.Item(1, 7).Value = "持有市值(亿元)"
$ws.Cells.Item(1, 8).Value = "仓位排名"

# Fund rows: (row, code, name, scale, stockPosition, positionRatio, marketValue, rank)
$data = @(
    @(2, "502000", "西部利得中证500指数增强（LOF）A", "30.03", "85.49", "1.98", "0.5946", 10),
    @(3, "006593", "博道中证500指数增强A", "31.17", "93.08", "1.24", "0.3865", 9),
    @(4, "004148", "圆信永丰多策略精选混合", "7.48", "89.36", "3.84", "0.2872", 8),
    @(5, "010779", "西部利得量化优选一年持有期混合A", "11.33", "87.57", "2.30", "0.2606", 5),
    @(6, "006594", "博道中证500指数增强C", "10.41", "93.08", "1.24", "0.1291", 9),
    @(7, "009300", "西部利得中证500指数增强（LOF）C", "6.33", "85.49", "1.98", "0.1253", 10),
    @(8, "007831", "博道伍佰智航股票A", "10.01", "92.13", "1.14", "0.1141", 3),
    @(9, "007126", "博道远航混合A", "10.18", "88.61", "1.05", "0.1069", 10),
    @(10, "202019", "南方策略优化混合", "3.33", "94.36", "2.50", "0.0832", 7),
    @(11, "005347", "诺德量化优选6个月持有期混合", "2.60", "93.66", "2.84", "0.0738", 9),
    @(12, "970041", "国海量化优选一年持有股票A", "7.70", "91.93", "0.95", "0.0732", 1),
    @(13, "970042", "国海量化优选一年持有股票C", "6.95", "91.93", "0.95", "0.0660", 1),
    @(14, "007832", "博道伍佰智航股票C", "5.65", "92.13", "1.14", "0.0644", 3),
    @(15, "007127", "博道远航混合C", "5.94", "88.61", "1.05", "0.0624", 10),
    @(16, "006267", "诺德量化核心灵活配置混合A", "1.84", "93.91", "2.87", "0.0528", 8),
    @(17, "006969", "圆信永丰高端制造混合", "1.04", "86.60", "4.03", "0.0419", 6),
    @(18, "010780", "西部利得量化优选一年持有期混合C", "1.44", "87.57", "2.30", "0.0331", 5),
    @(19, "000270", "建信灵活配置混合", "2.12", "93.93", "0.99", "0.0210", 6),
    @(20, "005381", "泰康睿利量化多策略混合A", "0.99", "93.49", "2.03", "0.0201", 5),
    @(21, "006268", "诺德量化核心灵活配置混合C", "0.50", "93.91", "2.87", "0.0144", 8),
    @(22, "005382", "泰康睿利量化多策略混合C", "0.49", "93.49", "2.03", "0.0099", 5),
    @(23, "590007", "中邮中证500指数增强A", "0.43", "91.51", "1.59", "0.0068", 5),
    @(24, "005120", "上投摩根量化多因子灵活配置混合", "0.21", "92.54", "2.71", "0.0057", 6),
    @(25, "003717", "中银量化精选灵活配置混合A", "0.49", "90.38", "1.15", "0.0056", 9),
    @(26, "008124", "中邮中证500指数增强C", "0.04", "91.51", "1.59", "0.0006", 5),
    @(27, "010484", "中银量化精选灵活配置混合C", "0.01", "90.38", "1.15", "0.0001", 9)
)

foreach ($row in $data) {
    $r = $row[0]

    # Column A: zero-based running index (0..25), numeric.
    $ws.Cells.Item($r, 1).Value = ($r - 2)

    # Columns B..G are stored as TEXT in the source workbook (even though
    # several look numeric, e.g. "30.03"), so force a text format before
    # assigning the values.
    $textRange = $ws.Range("B" + $r + ":G" + $r)
    $textRange.NumberFormat = "@"

    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $ws.Cells.Item($r, 7).Value = $row[6]

    # Column H: rank, numeric.
    $ws.Cells.Item($r, 8).Value = $row[7]
}

# ---------------------------------------------------------------------
# 2. Insert a new top data row into the "总计" (summary) sheet for the
#    freshly-added "2022-Q1" quarter, shifting the previous rows down.
# ---------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Re-use row 3's (formerly row 2's) style for the new index cell A2 so the
# bold/centered/bordered look carries over to the new row, matching the
# style already applied to every other index cell in column A.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 26
$total.Range("D2").Value = 2.64

# The pre-existing rows keep their relative order but their running index in
# column A shifts up by one (0,1,2,3 -> 1,2,3,4) now that a new row 0 exists.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3
$total.Range("A6").Value = 4
